$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of results
$ws.Range("A22").Value = "Bandpower + PCA + NuSVM (linear kernel)"
$ws.Range("B22").Value = 0.9605
$ws.Range("B22").NumberFormat = $ws.Range("B21").NumberFormat
$ws.Range("C22").Value = "19/19"
$ws.Range("D22").Value = "RH"
$ws.Range("E22").Value = "0, 1, 1, 2, 3, 3, 5, 12, 13, 23, 30, 52, 57"
$ws.Range("F22").Value = "nu=8585, n_components=3, freq bands (Hz) 4-8,8-13,13-30"

# Update view to reflect new active cell / scroll position
$ws.Range("F22").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 5
